# Regenerate the localization handoff/handback status report.
#
# The CI job that produces this workbook re-ran and a handful of rows whose
# "Ready for handoff" / "Handback transform failed" status hadn't actually
# progressed got re-stamped with the latest run's timestamp (the job also
# collapsed a couple of now-duplicate timestamp strings in the process).
# Net effect on the data: the "Latest Handoff Date" / "Latest Handoff
# Datetime" columns for the still-pending rows all take on the newest
# timestamp for that run.

$wb = $excel.ActiveWorkbook

$overviewRows = @(7, 10, 11, 12, 13, 14, 15, 16)
$langRows     = @(7, 10, 11, 12, 13, 14, 15, 16)

# "Overview" sheet - column D ("Latest Handoff Date")
$ws = $wb.Worksheets.Item("Overview")
foreach ($r in $overviewRows) {
    $ws.Range("D$r").Value = "2016-25-13 06:25:54"
}

# "zh-cn" sheet - column E ("Latest Handoff Datetime")
$ws = $wb.Worksheets.Item("zh-cn")
foreach ($r in $langRows) {
    $ws.Range("E$r").Value = "2016-03-13 06:25:50"
}

# "de-de" sheet - column E ("Latest Handoff Datetime")
$ws = $wb.Worksheets.Item("de-de")
foreach ($r in $langRows) {
    $ws.Range("E$r").Value = "2016-03-13 06:25:54"
}
